$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.950.47'
$ws.Range("E2").Value = '  +4.01%  '
$ws.Range("D3").Value = '2.780.32'
$ws.Range("E3").Value = '  +4.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.85%  '
$ws.Range("E7").Value = '  +3.72%  '
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").Value = '2.805.55'
$ws.Range("E9").Value = '  +4.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.24%  '
$ws.Range("E11").Value = '  +3.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.398'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.76%  '
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").Value = '3.278.60'
$ws.Range("E14").Value = '  +4.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.70%  '
$ws.Range("D16").Value = '63.888.45'
$ws.Range("E16").Value = '  +4.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000160'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +8.21%  '
$ws.Range("D18").Value = '2.792.48'
$ws.Range("E18").Value = '  +4.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '368.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.552'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.176'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.93%  '
$ws.Range("E27").Value = '  +2.61%  '
$ws.Range("D28").Value = '0.0₃0972'
$ws.Range("E28").Value = '  +16.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +1.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.76%  '
$ws.Range("E32").Value = '  +10.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '173.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.11'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.89'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.27%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.51'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.61%  '
$ws.Range("E38").Value = '  +7.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '344.49'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.29'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +15.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0613'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.29%  '
$ws.Range("E47").Value = '  +3.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0263'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("D51").Value = '2.189.07'
$ws.Range("E51").Value = '  +3.43%  '
